$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows to make room for the additional "book" detail rows.
# Original rows 5-10 become rows 6,9,10,12,13,14 after these inserts,
# and the newly inserted rows become 5,7,8,11.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(11).Insert()

# New row 5: Good Actions - Root Text Book
$ws.Range("E5").Value = "ཀུན་བཟང་སྨོན་ལམ་རྩ་བ་དེབ་།"
$ws.Range("F5").Value = "The Prayer of Good Actions Root Text Book"

# New row 7: Good Actions - Commentary Text Book 1
$ws.Range("E7").Value = "ཀུན་བཟང་སྨོན་ལམ་འགྲེལ་བ་དེབ་། ༡"
$ws.Range("F7").Value = "The Prayer of Good Actions Commentary Text Book 1"

# New row 8: Good Actions - Commentary Text Book 2
$ws.Range("E8").Value = "ཀུན་བཟང་སྨོན་ལམ་འགྲེལ་བ་དེབ་། ༢"
$ws.Range("F8").Value = "The Prayer of Good Actions Commentary Text Book 2"

# New row 11: Good Conduct - Root Text Book
$ws.Range("E11").Value = "བཟང་སྤྱོད་སྨོན་ལམ་རྩ་བ་དེབ་།"
$ws.Range("F11").Value = "The Prayer of Good Conduct Root Text Book"

# Existing rows 13/14 (formerly 9/10) gain Good Conduct commentary book details.
# Copy formatting from the neighboring already-styled cell first so the new
# E cells pick up the same cell style as the rest of the column.
$ws.Range("B13").Copy($ws.Range("E13"))
$ws.Range("E13").Value = "བཟང་སྤྱོད་སྨོན་ལམ་འགྲེལ་བ་དེབ་། ༡"
$ws.Range("F13").Value = "The Prayer of Good Conduct Commentary Text Book 1"

$ws.Range("B14").Copy($ws.Range("E14"))
$ws.Range("E14").Value = "བཟང་སྤྱོད་སྨོན་ལམ་འགྲེལ་བ་དེབ་། ༢"
$ws.Range("F14").Value = "The Prayer of Good Conduct Commentary Text Book 2"

# Size the new column E to match the authored width (stored width of 16).
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666

# Match the author's final cursor position when they saved the file.
$null = $ws.Range("E14").Select()
